# Generate Report for Handoff
# Refresh the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps
# for the rows that were re-handed-off since the last report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows (2-16) whose "Latest Handoff" timestamp needs to be refreshed.
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $rows) {
    $wsOverview.Range("D$r").Value = "2016-03-22 20:31:29"
    $wsZhCn.Range("E$r").Value     = "2016-03-22 20:31:24"
    $wsDeDe.Range("E$r").Value     = "2016-03-22 20:31:29"
}

$wb.Save()
